# Apply the commit:
#  1. Clear a handful of stray empty inline-string cells left over on
#     Table_1 (they carry no value, they're just noise from the source
#     export) so the row/cell set matches the cleaned-up sheet.
#  2. Add a second worksheet "Table_2" (right after "Table_1") holding
#     the capital-adequacy ratio summary table.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- 1. Drop the empty inlineStr cells on Table_1 -------------------------
$emptyCells = @("B2", "A3", "B19", "A20", "B21", "B28", "B30", "B41", "B42")
foreach ($addr in $emptyCells) {
    $ws1.Range($addr).ClearContents()
}

# --- 2. Insert the new "Table_2" worksheet right after "Table_1" ---------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Table_2"

# Header row (bold/centered/bordered style, matching the other header row
# used elsewhere in this workbook).
$ws2.Range("A1").Value = "Əmsal"
$ws2.Range("B1").Value = "Norma (Sistem əhəmiyyətli)"
$ws2.Range("C1").Value = "Norma (Banklar istisna)"
$ws2.Range("D1").Value = "Fakt"
$ws1.Range("A1:B1").Copy()
$ws2.Range("A1:D1").PasteSpecial(-4122)  # xlPasteFormats

# Data rows - keep the percentages / "minimum N%" figures as plain text.
$ws2.Range("A2:D4").NumberFormat = "@"

$ws2.Range("A2").Value = "9.  I dərəcəli  kapitalın  adekvatlıq əmsalı"
$ws2.Range("B2").Value = "6.0%"
$ws2.Range("C2").Value = "5.0%"
$ws2.Range("D2").Value = "9.6%"

$ws2.Range("A3").Value = "10. məcmu kapitalın  adekvatlıq  əmsalı"
$ws2.Range("B3").Value = "12.0%"
$ws2.Range("C3").Value = "10.0%"
$ws2.Range("D3").Value = "12.2%"

$ws2.Range("A4").Value = "11. Leverec əmsalı"
$ws2.Range("B4").Value = "minimum 5%"
$ws2.Range("C4").Value = "minimum 4%"
$ws2.Range("D4").Value = "5.9%"
